$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the "Good Morning" greeting cell to the new commit text
$ws.Range("E8").Value = "GIT UPDATE"

# Mirror the author's selection ending on the edited cell
$ws.Range("E8").Select()
